# Zotero3 edit: add "Another reference " citation paragraph + matching
# bibliography entry for Hume, Kathryn's "Attenuated Realities" article.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1) New paragraph after the first ("A reference ...") paragraph: a second
#        Zotero citation field "Another reference (Hume)". ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()

$para1Runs = '<w:r><w:t xml:space="preserve">Another reference </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> ADDIN ZOTERO_ITEM CSL_CITATION {"citationID":"l1UNDPDQ","properties":{"formattedCitation":"(Hume)","plainCitation":"(Hume)"},"citationItems":[{"id":8,"uris":["http://zotero.org/users/local/9yJPt4JR/items/XVNU7KZU"],"uri":["http://zotero.org/users/local/9yJPt4JR/items/XVNU7KZU"],"itemData":{"id":8,"type":"article-journal","title":"Attenuated Realities: Pynchon''s Trajectory from V. to Inherent Vice","container-title":"Orbit: Writing Around Pynchon","volume":"2","issue":"1","source":"www.pynchon.net","abstract":"Part of what established Pynchon as postmodern was his piling up of multiple realities. Hence, the surprise that  Inherent Vice  retains only the most attenuated forms of such worlds alternative to our own. In earlier fiction, we find a world served by the Tristero postal system, a world inhabited by angels, by thanatoids and other forms of the dead, by Japanese film monsters, by giant vegetables. In Pynchon’s fictive realities, an airship can sail beneath desert sand or through the center of the globe via Symmes’s Hole, and the photograph of a corpse can be run backward in time to show its murderer. Up through  Against the Day , Pynchon showered us with alternate realities that reached beyond the material world that most of us accept as  alles, was der Fall ist .  Inherent Vice  departs from this vision. Has Pynchon simply grown up? Or grown old? Or is something else operating here? I will provide a brief taxonomy of Pynchon’s multiple worlds as characterized by paranoia, mysticism, religion, and humor and then analyze what remains of these in  Inherent Vice . Among the causes for his changed technique may be his choice of genre. The detective story is epistemological rather than ontological in its questions, so Pynchon concerns himself far more with what Doc Sportello can know than with making him navigate through multiple realities. I argue, however, that  Inherent Vice  is surprisingly a worst-case scenario for Pynchon.","URL":"https://www.pynchon.net/owap/article/view/50","DOI":"10.7766/orbit.v2.1.50","ISSN":"2047-2870","shortTitle":"Attenuated Realities","language":"en","author":[{"family":"Hume","given":"Kathryn"}],"issued":{"date-parts":[["2013",12,12]]},"accessed":{"date-parts":[["2014",1,15]],"season":"12:02:44"}}}],"schema":"https://github.com/citation-style-language/schema/raw/master/csl-citation.json"} </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>(Hume)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r>'

$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.Collapse(1)
$r2.InsertXML($pkgHeader + '<w:p>' + $para1Runs + '</w:p>' + $pkgFooter)

# --- 2) New Bibliography-styled paragraph after the existing Herman entry,
#        containing the Hume, Kathryn bibliography entry. ---
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()

$bibParaRuns = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve">Hume, Kathryn. “Attenuated Realities: Pynchon’s Trajectory from V. to Inherent Vice.” </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:i/><w:iCs/></w:rPr><w:t>Orbit: Writing Around Pynchon</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve"> 2.1 (2013): n. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>pag</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:i/><w:iCs/></w:rPr><w:t>www.pynchon.net</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>. Web. 15 Jan. 2014.</w:t></w:r>'

$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.Collapse(1)
$bibPPr = '<w:pPr><w:pStyle w:val="Bibliography"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr>'
$r4.InsertXML($pkgHeader + '<w:p>' + $bibPPr + $bibParaRuns + '</w:p>' + $pkgFooter)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    Write-Host "[$i] $($p.Range.Text)"
}
